$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1037.75
$ws.Range("I19").Value = 850
$ws.Range("J19").Value = 1100.3334
$ws.Range("K19").Value = 850
$ws.Range("L19").Value = 1100.3334
$ws.Range("M19").Value = -675
$ws.Range("N19").Value = -1450.3334
$ws.Range("H40").Value = 6888.231
$ws.Range("I40").Value = 5650.143
$ws.Range("K40").Value = 5650.143
$ws.Range("M40").Value = -5475.143
$ws.Range("H53").Value = 399.7143
$ws.Range("I53").Value = 520.125
$ws.Range("J53").Value = 239.16667
$ws.Range("K53").Value = 520.125
$ws.Range("L53").Value = 239.16667
$ws.Range("M53").Value = 116.875
$ws.Range("N53").Value = -1513.16667
$ws.Range("H70").Value = 3525.85
$ws.Range("I70").Value = 2200.4285
$ws.Range("J70").Value = 4239.5386
$ws.Range("K70").Value = 6601.2855
$ws.Range("L70").Value = 12718.6158
$ws.Range("M70").Value = -6331.2855
$ws.Range("N70").Value = -13258.6158
$ws.Range("H73").Value = 3525.85
$ws.Range("I73").Value = 2200.4285
$ws.Range("J73").Value = 4239.5386
$ws.Range("K73").Value = 6601.2855
$ws.Range("L73").Value = 12718.6158
$ws.Range("M73").Value = -5665.2855
$ws.Range("N73").Value = -14590.6158
$ws.Range("H88").Value = 1913.1666
$ws.Range("I88").Value = 1950
$ws.Range("J88").Value = 1894.75
$ws.Range("K88").Value = 1950
$ws.Range("L88").Value = 1894.75
$ws.Range("M88").Value = -1544
$ws.Range("N88").Value = -2706.75
$ws.Range("H91").Value = 1913.1666
$ws.Range("I91").Value = 1950
$ws.Range("J91").Value = 1894.75
$ws.Range("K91").Value = 1950
$ws.Range("L91").Value = 1894.75
$ws.Range("M91").Value = -546
$ws.Range("N91").Value = -4702.75
$ws.Range("H100").Value = 1606.7778
$ws.Range("I100").Value = 1370.125
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 1370.125
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -829.125
$ws.Range("N100").Value = -4582
$ws.Range("H101").Value = 542
$ws.Range("I101").Value = 469
$ws.Range("J101").Value = 615
$ws.Range("K101").Value = 1407
$ws.Range("L101").Value = 1845
$ws.Range("M101").Value = 215
$ws.Range("N101").Value = -5089

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 575.5
$ws.Range("J4").Value = 1111
$ws.Range("L4").Value = 1111
$ws.Range("N4").Value = -1343
$ws.Range("H5").Value = 122.954544
$ws.Range("J5").Value = 119
$ws.Range("L5").Value = 119
$ws.Range("N5").Value = -343

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 122.954544
$ws.Range("J4").Value = 119
$ws.Range("L4").Value = 119
$ws.Range("N4").Value = -349
$ws.Range("H58").Value = 150000
$ws.Range("J58").Value = 150000
$ws.Range("L58").Value = 150000
$ws.Range("N58").Value = -150588
$ws.Range("H105").Value = 4133535.5
$ws.Range("I105").Value = 5348671
$ws.Range("K105").Value = 5348671
$ws.Range("M105").Value = -5346924

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H99").Value = 3666.6667
$ws.Range("I99").Value = 3666.6667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3666.6667
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 3666.6667
$ws.Range("I126").Value = 3666.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11000.0001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1788.1904
$ws.Range("I132").Value = 1788.1904
$ws.Range("K132").Value = 5364.5712
$ws.Range("M132").Value = -2834.5712

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 25151
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 200.6875
$ws.Range("I12").Value = 107.5
$ws.Range("K12").Value = 322.5
$ws.Range("M12").Value = -149.5
$ws.Range("H26").Value = 171
$ws.Range("I26").Value = 171
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 513
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H61").Value = 175
$ws.Range("I61").Value = 150
$ws.Range("J61").Value = 225
$ws.Range("K61").Value = 450
$ws.Range("L61").Value = 675
$ws.Range("M61").Value = -235
$ws.Range("N61").Value = -1105

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 24091428
$ws.Range("I7").Value = 17797778
$ws.Range("K7").Value = 17797778
$ws.Range("M7").Value = -17797666
$ws.Range("H8").Value = 24091428
$ws.Range("I8").Value = 17797778
$ws.Range("K8").Value = 17797778
$ws.Range("M8").Value = -17797639
$ws.Range("H19").Value = 1225
$ws.Range("I19").Value = 1200
$ws.Range("J19").Value = 1250
$ws.Range("K19").Value = 1200
$ws.Range("L19").Value = 1250
$ws.Range("M19").Value = -912
$ws.Range("N19").Value = -1826

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1291.6666
$ws.Range("I7").Value = 1300
$ws.Range("K7").Value = 1300
$ws.Range("M7").Value = -1188
$ws.Range("H11").Value = 16333.333
$ws.Range("J11").Value = 16333.333
$ws.Range("L11").Value = 16333.333
$ws.Range("N11").Value = -16613.333
$ws.Range("H22").Value = 1222.3846
$ws.Range("I22").Value = 898.6667
$ws.Range("J22").Value = 1499.8572
$ws.Range("K22").Value = 898.6667
$ws.Range("L22").Value = 1499.8572
$ws.Range("M22").Value = -603.6667
$ws.Range("N22").Value = -2089.8572
$ws.Range("H27").Value = 1222.3846
$ws.Range("I27").Value = 898.6667
$ws.Range("J27").Value = 1499.8572
$ws.Range("K27").Value = 898.6667
$ws.Range("L27").Value = 1499.8572
$ws.Range("M27").Value = -791.6667
$ws.Range("N27").Value = -1713.8572
$ws.Range("H46").Value = 6972.278
$ws.Range("I46").Value = 5928.5713
$ws.Range("J46").Value = 7636.4546
$ws.Range("K46").Value = 5928.5713
$ws.Range("L46").Value = 7636.4546
$ws.Range("M46").Value = -5740.5713
$ws.Range("N46").Value = -8012.4546
$ws.Range("H61").Value = 76927880
$ws.Range("I61").Value = 111115040
$ws.Range("J61").Value = 6749.75
$ws.Range("K61").Value = 111115040
$ws.Range("L61").Value = 6749.75
$ws.Range("M61").Value = -111114838
$ws.Range("N61").Value = -7153.75
$ws.Range("H68").Value = 6443.75
$ws.Range("I68").Value = 1899
$ws.Range("K68").Value = 1899
$ws.Range("M68").Value = -1150
$ws.Range("H71").Value = 6443.75
$ws.Range("I71").Value = 1899
$ws.Range("K71").Value = 9495
$ws.Range("M71").Value = -5751
$ws.Range("H82").Value = 2891.3125
$ws.Range("I82").Value = 601.8570999999999
$ws.Range("J82").Value = 4672
$ws.Range("K82").Value = 601.8570999999999
$ws.Range("L82").Value = 4672
$ws.Range("M82").Value = -240.8570999999999
$ws.Range("N82").Value = -5394
$ws.Range("H85").Value = 2891.3125
$ws.Range("I85").Value = 601.8570999999999
$ws.Range("J85").Value = 4672
$ws.Range("K85").Value = 601.8570999999999
$ws.Range("L85").Value = 4672
$ws.Range("M85").Value = 646.1429000000001
$ws.Range("N85").Value = -7168
$ws.Range("H93").Value = 2247.5833
$ws.Range("I93").Value = 2297.3635
$ws.Range("J93").Value = 1700
$ws.Range("K93").Value = 2297.3635
$ws.Range("L93").Value = 1700
$ws.Range("M93").Value = -1049.3635
$ws.Range("N93").Value = -4196
$ws.Range("H100").Value = 7399.9
$ws.Range("I100").Value = 4799.8
$ws.Range("K100").Value = 4799.8
$ws.Range("M100").Value = -4258.8
$ws.Range("H113").Value = 76927880
$ws.Range("I113").Value = 111115040
$ws.Range("J113").Value = 6749.75
$ws.Range("K113").Value = 111115040
$ws.Range("L113").Value = 6749.75
$ws.Range("M113").Value = -111112870
$ws.Range("N113").Value = -11089.75
$ws.Range("H122").Value = 6186.5
$ws.Range("I122").Value = 5998.8335
$ws.Range("J122").Value = 6749.5
$ws.Range("K122").Value = 17996.5005
$ws.Range("L122").Value = 20248.5
$ws.Range("M122").Value = -15546.5005
$ws.Range("N122").Value = -25148.5
$ws.Range("H126").Value = 1291.6666
$ws.Range("I126").Value = 1300
$ws.Range("K126").Value = 3900
$ws.Range("M126").Value = -1430

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 687.8
$ws.Range("I100").Value = 734.75
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 1469.5
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -928.5
$ws.Range("N100").Value = -2082
$ws.Range("H107").Value = 41667536
$ws.Range("J107").Value = 1739.5
$ws.Range("L107").Value = 5218.5
$ws.Range("N107").Value = -9058.5
$ws.Range("H132").Value = 2922.75
$ws.Range("I132").Value = 2846.75
$ws.Range("K132").Value = 8540.25
$ws.Range("M132").Value = -6010.25
